# setGameModal fix game time
# Insert a new game ("Resident Evil 5 ") into the sorted game-time table.
# The table is sorted descending by column D (Время, ч / Time, hours),
# so the new row is inserted at its correct sorted position (row 52),
# pushing the former rows 52-62 down to 53-63.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Make room for the new row by inserting a blank row at position 52.
$ws.Rows.Item(52).Insert()

# Fill in the new game's data.
$ws.Cells.Item(52, 1).Value = "Resident Evil 5 "
$ws.Cells.Item(52, 3).Value = 9
$ws.Cells.Item(52, 4).Value = 6.3333333333333304

# Update the view to match where the edit happened.
$ws.Range("D51").Select()
$excel.ActiveWindow.ScrollRow = 31

$wsView = $ws.Application.ActiveWindow
$ws.Range("A31").Select()
$ws.Range("D51").Select()
